$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.436.21'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '1.551.83'
$ws.Range('E3').Value = '  -1.89%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.483'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.81%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  +1.18%  '
$ws.Range('E9').Value = '  -1.96%  '
$ws.Range('E10').Value = '  -1.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0891'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('D12').Value = '1.776.34'
$ws.Range('E12').Value = '  -1.72%  '
$ws.Range('D13').Value = '1.562.46'
$ws.Range('E13').Value = '  -0.07%  '
$ws.Range('D14').Value = '28.442.15'
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.63'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.509'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.02'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.03'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.69%  '
$ws.Range('D20').Value = '0.0₃0673'
$ws.Range('E20').Value = '  -2.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E22').Value = '  -1.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.90'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.03'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.12'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('E26').Value = '  -1.88%  '
$ws.Range('E27').Value = '  -1.34%  '
$ws.Range('E28').Value = '  -0.13%  '
$ws.Range('E29').Value = '  -2.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0465'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.88%  '
$ws.Range('E31').Value = '  -4.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.15'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.73%  '
$ws.Range('D33').Value = '1.384.84'
$ws.Range('E33').Value = '  -0.91%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.99'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.04'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.53%  '
$ws.Range('E36').Value = '  -1.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.29'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.65'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0162'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.86%  '
$ws.Range('E40').Value = '  +2.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.512'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.79%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.770'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.43%  '
$ws.Range('E44').Value = '  -0.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.34'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.74'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('D47').Value = '1.688.43'
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.869'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.97%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '85.27'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '42.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.21%  '
$ws.Range('D51').Value = '0.0₆0101'
$ws.Range('E51').Value = '  -1.65%  '
